# Update the "想去人数" (F column) figures in both the "展览" and
# "全部类型" worksheets to match the refreshed data snapshot.
#
# Mapping of row -> new F-column value (same on both sheets):
#   F8:  33   -> 34
#   F11: 117  -> 118
#   F12: 1142 -> 1147
#   F13: 1468 -> 1472
#   F17: 101  -> 104
#   F20: 98   -> 99
#   F23: 303  -> 304
#   F24: 1679 -> 1681
#   F28: 633  -> 635
#   F30: 107  -> 116
#   F31: 3981 -> 3989
#   F32: 10   -> 11
#   F34: 239  -> 240
#   F35: 1017 -> 1022
#   F36: 107  -> 108
#   F39: 116  -> 119

$wb = $excel.ActiveWorkbook

$updates = @{
    8  = 34
    11 = 118
    12 = 1147
    13 = 1472
    17 = 104
    20 = 99
    23 = 304
    24 = 1681
    28 = 635
    30 = 116
    31 = 3989
    32 = 11
    34 = 240
    35 = 1022
    36 = 108
    39 = 119
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
